# Apply updated financial figures across the CFO Dashboard workbook.
$wb = $excel.ActiveWorkbook

# --- LIQUIDITY_MONITOR ---
$wsLiquidity = $wb.Worksheets.Item("LIQUIDITY_MONITOR")
$wsLiquidity.Range("B5").Value = 56421
$wsLiquidity.Range("B6").Value = 28013
$wsLiquidity.Range("B9").Value = 56421

$wsLiquidity.Range("C19").Value = 295885
$wsLiquidity.Range("E19").Value = 0
$wsLiquidity.Range("G19").Value = 0
$wsLiquidity.Range("I19").Value = 0

# E19/G19/I19 switch from the "HARD value" shading (style of C19) to the
# plain zero shading already used by D19/F19/H19 on the same row.
$wsLiquidity.Range("D19").Copy()
$wsLiquidity.Range("E19").PasteSpecial(-4122)
$wsLiquidity.Range("D19").Copy()
$wsLiquidity.Range("G19").PasteSpecial(-4122)
$wsLiquidity.Range("D19").Copy()
$wsLiquidity.Range("I19").PasteSpecial(-4122)

$wsLiquidity.Range("C20").Value = -17468
$wsLiquidity.Range("D20").Value = -61630
$wsLiquidity.Range("F20").Value = -11620
$wsLiquidity.Range("G20").Value = -53250

# --- PROFIT_CONTROL ---
$wsProfit = $wb.Worksheets.Item("PROFIT_CONTROL")
$wsProfit.Range("B5").Value = 0.593221527602339
$wsProfit.Range("B6").Value = 0.04422913950594023

$wsProfit.Range("B11").Value = 1183541
$wsProfit.Range("C11").Value = 1183541
$wsProfit.Range("B12").Value = 481439
$wsProfit.Range("B13").Value = 702102
$wsProfit.Range("B15").Value = 361875
$wsProfit.Range("C15").Value = 43907
$wsProfit.Range("B17").Value = 52347

# --- BALANCE_SHEET_HEALTH ---
$wsBalance = $wb.Worksheets.Item("BALANCE_SHEET_HEALTH")
$wsBalance.Range("B8").Value = 183282

# --- DEBT_MANAGER ---
$wsDebt = $wb.Worksheets.Item("DEBT_MANAGER")
$wsDebt.Range("B14").Value = 183282
